# Apply cryptos.xlsx data refresh (prices & 1h volume changes) as of
# Tue May 30 20:09:29 UTC 2023, plus a row-order swap for rows 39/40.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    # Use Value2 + a leading apostrophe to force text storage even when the
    # string looks numeric (e.g. "1.001", "0.5228"), matching the original
    # inline-string cell type. Reset style to Normal afterwards so no stray
    # "Text" number-format style gets attached to the cell (cells in this
    # sheet carry no explicit style).
    $ws.Range($cellRef).Value2 = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "27.874.44"
Set-TextCell "E2" "  +0.66%  "
Set-TextCell "D3" "1.908.86"
Set-TextCell "E3" "  +0.74%  "
Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  -0.11%  "
Set-TextCell "D5" "313.35"
Set-TextCell "E5" "  +0.35%  "
Set-TextCell "E6" "  -0.03%  "
Set-TextCell "D7" "0.5228"
Set-TextCell "E7" "  +7.49%  "
Set-TextCell "D8" "0.3797"
Set-TextCell "E8" "  +0.02%  "
Set-TextCell "D9" "0.07246"
Set-TextCell "E9" "  -1.04%  "
Set-TextCell "D10" "21.38"
Set-TextCell "E10" "  +4.07%  "
Set-TextCell "D11" "0.9104"
Set-TextCell "E11" "  -0.37%  "
Set-TextCell "D12" "1.939.91"
Set-TextCell "E12" "  +2.60%  "
Set-TextCell "D13" "0.07649"
Set-TextCell "E13" "  -0.22%  "
Set-TextCell "D14" "5.460"
Set-TextCell "E14" "  -0.18%  "
Set-TextCell "D15" "92.40"
Set-TextCell "E15" "  +1.40%  "
Set-TextCell "D16" "1.001"
Set-TextCell "E16" "  -0.13%  "
Set-TextCell "D17" "0.000008716"
Set-TextCell "D18" "1.001"
Set-TextCell "D19" "27.916.04"
Set-TextCell "E19" "  +0.99%  "
Set-TextCell "D20" "14.57"
Set-TextCell "E20" "  +0.48%  "
Set-TextCell "E21" "  +0.69%  "
Set-TextCell "D22" "2.163.36"
Set-TextCell "E22" "  +2.44%  "
Set-TextCell "D23" "10.87"
Set-TextCell "E23" "  +1.09%  "
Set-TextCell "D24" "6.645"
Set-TextCell "E24" "  +0.61%  "
Set-TextCell "D25" "153.70"
Set-TextCell "E25" "  -0.01%  "
Set-TextCell "D26" "1.871"
Set-TextCell "E26" "  -1.86%  "
Set-TextCell "D27" "2.173"
Set-TextCell "E27" "  +1.62%  "
Set-TextCell "D28" "18.35"
Set-TextCell "E28" "  -0.15%  "
Set-TextCell "D29" "114.83"
Set-TextCell "E29" "  -0.58%  "
Set-TextCell "D30" "4.878"
Set-TextCell "E30" "  -0.06%  "
Set-TextCell "D31" "0.09015"
Set-TextCell "E31" "  +1.38%  "
Set-TextCell "D32" "4.876"
Set-TextCell "E32" "  +5.21%  "
Set-TextCell "D33" "3.184"
Set-TextCell "E33" "  -0.57%  "
Set-TextCell "D34" "1.235"
Set-TextCell "E34" "  +1.04%  "
Set-TextCell "D35" "0.7815"
Set-TextCell "E35" "  +1.91%  "
Set-TextCell "D36" "0.02102"
Set-TextCell "E36" "  +3.33%  "
Set-TextCell "D37" "2.617"
Set-TextCell "E37" "  +3.34%  "
Set-TextCell "D38" "3.078"
Set-TextCell "E38" "  +3.36%  "
Set-TextCell "B39" "TheSandbox"
Set-TextCell "C39" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D39" "0.5588"
Set-TextCell "E39" "  +2.06%  "
Set-TextCell "B40" "TrustWalletToken"
Set-TextCell "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D40" "1.094"
Set-TextCell "E40" "  -0.05%  "
Set-TextCell "D41" "0.05292"
Set-TextCell "E41" "  +0.38%  "
Set-TextCell "D42" "6.725"
Set-TextCell "E42" "  -2.55%  "
Set-TextCell "D43" "115.54"
Set-TextCell "E43" "  +3.49%  "
Set-TextCell "D44" "8.579"
Set-TextCell "E44" "  +0.98%  "
Set-TextCell "D45" "0.1518"
Set-TextCell "E45" "  -0.17%  "
Set-TextCell "D46" "0.4821"
Set-TextCell "E46" "  +0.66%  "
Set-TextCell "D47" "10.45"
Set-TextCell "E47" "  -1.49%  "
Set-TextCell "E48" "  +0.02%  "
Set-TextCell "D49" "1.625"
Set-TextCell "E49" "  -0.40%  "
Set-TextCell "D50" "67.09"
Set-TextCell "E50" "  -0.35%  "
Set-TextCell "D51" "0.05998"
Set-TextCell "E51" "  -0.85%  "
